$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: ALINE BASIOTE entry removed -> table shifts up by one.
# Row 19 becomes what used to be row 20 (ANTONIO VIEIRA), but its date moves
# from 44851 to 44852.
$ws.Range("B19").Value = "ANTONIO VIEIRA"
$ws.Range("C19").Value = "533ae974ff5ab9eaed4a9f8074909ec5"
$ws.Range("D19").Value = 44852
$ws.Range("G19").Value = "VENDA 16 (18/10)"

# --- Row 20: becomes what used to be row 21 (BEBETO SANTOS) - values unchanged.
$ws.Range("B20").Value = "BEBETO SANTOS"
$ws.Range("C20").Value = "5f4cb48c0ad331085484acd234d7f73b"
$ws.Range("D20").Value = 44852
$ws.Range("G20").Value = "VENDA 17 (18/10)"

# --- Row 21: becomes what used to be row 22 (EDENILSON SILVA) - values unchanged.
$ws.Range("B21").Value = "EDENILSON SILVA"
$ws.Range("C21").Value = "e26d1cd1918b4c7d99e4509543ea983a"
$ws.Range("D21").Value = 44852
$ws.Range("G21").Value = "VENDA 18 (18/10)"

# --- Rows for CESAR AUGUSTO / FERNANDO BRITO / DANIELE LOPES are deleted
# entirely, so the remaining records compact further upward.

# --- Row 22: becomes SERGIO (SOCIO JARDSON), date pushed out to 44855, and the
# previously-blank venda cell now holds the "-" placeholder.
$ws.Range("B22").Value = "SERGIO (SOCIO JARDSON)"
$ws.Range("C22").Value = "20f89519cd4d6fd819360d29f87f0df1"
$ws.Range("D22").Value = 44855
$ws.Range("G22").Value = "-"

# --- Row 23: becomes LUCIO GO.
$ws.Range("B23").Value = "LUCIO GO"
$ws.Range("C23").Value = "6d2ba0b66571df1a0f30a4ba316c2df8"
$ws.Range("D23").Value = 44856
$ws.Range("G23").Value = "VENDA 22 (22/10)"

# --- Row 24: becomes RENATO (SOCIO JARDSON), venda cell now holds "-".
$ws.Range("B24").Value = "RENATO (SOCIO JARDSON)"
$ws.Range("C24").Value = "3789612e62b0a636b6149d0d5cfbfc79"
$ws.Range("D24").Value = 44856
$ws.Range("G24").Value = "-"

# --- Row 25: becomes GABRIEL DE DEUS.
$ws.Range("B25").Value = "GABRIEL DE DEUS"
$ws.Range("C25").Value = "84e7824334195ec0675c3a0bde9b8bf4"
$ws.Range("D25").Value = 44863
$ws.Range("G25").Value = "VENDA 22 (29/10)"

# --- Row 26: brand-new client record appended at the end of the list.
$ws.Range("B26").Value = "MARIO FELIPE PEDROZO"
$ws.Range("C26").Value = "da64cb3429de3cb92e83d5cb9e2d8f9e"
$ws.Range("D26").Value = 44868
$ws.Range("G26").Value = "VENDA 23 (03/11)"

# --- Rows 27-29 no longer hold any client records - clear them out, copying
# the blank template row's date-cell number format (row 30) so the empty
# date cell formatting matches the rest of the unused rows below.
$ws.Range("D30").Copy()
$ws.Range("D27:D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B27:G29").ClearContents()

# --- Update the saved view state: zoomed to 75%, scrolled back to the top,
# with B27 as the active selection.
$excel.ActiveWindow.Zoom = 75
$ws.Range("B27").Select()
